# PlanTest-Script2.xlsx edit: add the JSON `storeKeys(json,jsonpath,var)`
# command to the hidden "#system" lookup sheet's `json` column, and remove
# the `text` entry from the `target` column (the "text" named-range's
# backing column (Y) is also being dropped entirely, shifting the
# web/webalert/webcookie/ws/ws.async/xml columns one to the left).
#
# NOTE: this runtime's Range.Insert()/Range.Delete() shift the WHOLE row
# (every column), not just the cells in the target range's columns, so a
# single-column insert/delete is instead done by hand: read/write the
# cell values directly, column by column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1) json (column M): insert a new row between the existing
#     storeValue/storeValues rows (M16/M17) and give it the new function
#     name. storeValue moves M16->M17, storeValues moves M17->M18.
$oldM17 = $ws.Cells.Item(17, 13).Value2
$ws.Cells.Item(18, 13).Value = $oldM17
$oldM16 = $ws.Cells.Item(16, 13).Value2
$ws.Cells.Item(17, 13).Value = $oldM16
$ws.Cells.Item(16, 13).Value = "storeKeys(json,jsonpath,var)"

# --- 2) target (column A): drop the "text" entry (was A25), shifting
#     web/webalert/webcookie/ws/ws.async/xml up one row (A26:A31 -> A25:A30).
for ($r = 25; $r -le 30; $r++) {
    $next = $ws.Cells.Item($r + 1, 1).Value2
    $ws.Cells.Item($r, 1).Value = $next
}
$ws.Cells.Item(31, 1).ClearContents()

# --- 3) text's backing column (Y) is removed outright, shifting the
#     web/webalert/webcookie/ws/ws.async/xml columns (Z:AE) one column left
#     (Z:AE -> Y:AD). This is a true full-column operation so it's safe to
#     use Columns(...).Delete() here.
$ws.Columns("Y").Delete()

# --- 4) Update the defined names so they reference the new ranges.
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
